$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-07 Tuesday", "2023-11-08 Wednesday"),
    @("37×83=", "71×19="),
    @("73×30=", "75×74="),
    @("96×74=", "69×52="),
    @("14×15=", "92×40="),
    @("56×19=", "35×24="),
    @("84×41=", "92×67="),
    @("48×20=", "97×31="),
    @("85×32=", "96×84="),
    @("30×32=", "89×20="),
    @("48×53=", "50×54="),
    @("71×18=", "67×12="),
    @("37×64=", "83×85="),
    @("36×99=", "22×17="),
    @("95×44=", "25×58="),
    @("74×23=", "68×56="),
    @("66×84=", "64×39="),
    @("36×64=", "30×76="),
    @("60×11=", "15×38="),
    @("42×68=", "25×71="),
    @("55×12=", "47×84="),
    @("93×79=", "44×71="),
    @("79×48=", "35×88="),
    @("91×12=", "24×61="),
    @("43×12=", "30×54="),
    @("57×23=", "67×32=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
